$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fbn1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 7.848425
$ws.Range("H2").Value = 23.545275
$ws.Range("I2").Value = 0.02436729568045431
$ws.Range("J2").Value = 0.02436729568045431
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 883.0529186076
$ws.Range("R2").Value = 7947.4762674684
$ws.Range("S2").Value = 0.007981051976556894
$ws.Range("T2").Value = 0.007981051976556894

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fbn1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 7.848425
$ws.Range("H3").Value = 23.545275
$ws.Range("I3").Value = 0.02436729568045431
$ws.Range("J3").Value = 0.02436729568045431
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 834.4011128160499
$ws.Range("R3").Value = 7509.61001534445
$ws.Range("S3").Value = 0.007541335870541443
$ws.Range("T3").Value = 0.007541335870541444

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fbn1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 7.848425
$ws.Range("H4").Value = 23.545275
$ws.Range("I4").Value = 0.02436729568045431
$ws.Range("J4").Value = 0.02436729568045431
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 978.6331050095167
$ws.Range("R4").Value = 8807.69794508565
$ws.Range("S4").Value = 0.00884490783335597
$ws.Range("T4").Value = 0.008844907833355972

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fbn1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 291.329961
$ws.Range("H5").Value = 873.989883
$ws.Range("I5").Value = 0.9045029162236017
$ws.Range("J5").Value = 0.9045029162236017
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 32778.52210333771
$ws.Range("R5").Value = 295006.6989300394
$ws.Range("S5").Value = 0.2962530139574873
$ws.Range("T5").Value = 0.2962530139574873

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fbn1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 291.329961
$ws.Range("H6").Value = 873.989883
$ws.Range("I6").Value = 0.9045029162236017
$ws.Range("J6").Value = 0.9045029162236017
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 30972.58923351582
$ws.Range("R6").Value = 278753.3031016424
$ws.Range("S6").Value = 0.2799309523952563
$ws.Range("T6").Value = 0.2799309523952563

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fbn1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 291.329961
$ws.Range("H7").Value = 873.989883
$ws.Range("I7").Value = 0.9045029162236017
$ws.Range("J7").Value = 0.9045029162236017
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 36326.41508528544
$ws.Range("R7").Value = 326937.735767569
$ws.Range("S7").Value = 0.328318949870858
$ws.Range("T7").Value = 0.3283189498708581

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fbn1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 22.91008466666667
$ws.Range("H8").Value = 68.730254
$ws.Range("I8").Value = 0.07112978809594397
$ws.Range("J8").Value = 0.07112978809594397
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 2577.691336853856
$ws.Range("R8").Value = 23199.2220316847
$ws.Range("S8").Value = 0.02329723180281213
$ws.Range("T8").Value = 0.02329723180281213

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fbn1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 22.91008466666667
$ws.Range("H9").Value = 68.730254
$ws.Range("I9").Value = 0.07112978809594397
$ws.Range("J9").Value = 0.07112978809594397
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 2435.673417351455
$ws.Range("R9").Value = 21921.06075616309
$ws.Range("S9").Value = 0.02201367067836857
$ws.Range("T9").Value = 0.02201367067836857

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fbn1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 22.91008466666667
$ws.Range("H10").Value = 68.730254
$ws.Range("I10").Value = 0.07112978809594397
$ws.Range("J10").Value = 0.07112978809594397
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 2856.696380913485
$ws.Range("R10").Value = 25710.26742822137
$ws.Range("S10").Value = 0.02581888561476328
$ws.Range("T10").Value = 0.02581888561476328
